# Remove the stray date stamp that was left in cell C1 of the "About" sheet.
# This drops both the cell's value (44307) and its date-formatted style,
# which in turn leaves that number-format style unused in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("C1").Clear()
